$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value2 = 44316
$ws.Cells.Item(2, 12).Value2 = 'Especial'
$ws.Cells.Item(2, 13).Value2 = 300
$ws.Cells.Item(2, 14).Value2 = 20000
$ws.Cells.Item(2, 15).Value2 = 20000
$ws.Cells.Item(2, 16).Value2 = 20000
$ws.Cells.Item(2, 17).Value2 = '$/caja 18 kilos granel'
$ws.Cells.Item(2, 18).Value2 = 'Provincia de Limarí'
$ws.Cells.Item(2, 19).Value2 = 1111
$ws.Cells.Item(2, 20).Value2 = 18

# Row 3
$ws.Cells.Item(3, 4).Value2 = 44299
$ws.Cells.Item(3, 12).Value2 = 'Primera'
$ws.Cells.Item(3, 13).Value2 = 100
$ws.Cells.Item(3, 14).Value2 = 15000
$ws.Cells.Item(3, 15).Value2 = 15000
$ws.Cells.Item(3, 16).Value2 = 15000
$ws.Cells.Item(3, 17).Value2 = '$/caja 15 kilos granel'
$ws.Cells.Item(3, 18).Value2 = 'Provincia de Curicó'
$ws.Cells.Item(3, 19).Value2 = 1000
$ws.Cells.Item(3, 20).Value2 = 15

# Row 4
$ws.Cells.Item(4, 4).Value2 = 44680
$ws.Cells.Item(4, 12).Value2 = 'Primera'
$ws.Cells.Item(4, 13).Value2 = 200
$ws.Cells.Item(4, 14).Value2 = 15000
$ws.Cells.Item(4, 15).Value2 = 15000
$ws.Cells.Item(4, 16).Value2 = 15000
$ws.Cells.Item(4, 17).Value2 = '$/caja 15 kilos granel'
$ws.Cells.Item(4, 18).Value2 = 'Provincia de Limarí'
$ws.Cells.Item(4, 19).Value2 = 1000
$ws.Cells.Item(4, 20).Value2 = 15

# Row 5
$ws.Cells.Item(5, 4).Value2 = 44291
$ws.Cells.Item(5, 12).Value2 = 'Primera'
$ws.Cells.Item(5, 13).Value2 = 150
$ws.Cells.Item(5, 14).Value2 = 12000
$ws.Cells.Item(5, 15).Value2 = 12000
$ws.Cells.Item(5, 16).Value2 = 12000
$ws.Cells.Item(5, 17).Value2 = '$/caja 15 kilos granel'
$ws.Cells.Item(5, 18).Value2 = 'Región Metropolitana'
$ws.Cells.Item(5, 19).Value2 = 800
$ws.Cells.Item(5, 20).Value2 = 15

# Row 6
$ws.Cells.Item(6, 4).Value2 = 44342
$ws.Cells.Item(6, 12).Value2 = 'Especial'
$ws.Cells.Item(6, 13).Value2 = 300
$ws.Cells.Item(6, 14).Value2 = 20000
$ws.Cells.Item(6, 15).Value2 = 20000
$ws.Cells.Item(6, 16).Value2 = 20000
$ws.Cells.Item(6, 17).Value2 = '$/caja 18 kilos granel'
$ws.Cells.Item(6, 18).Value2 = 'Provincia de Limarí'
$ws.Cells.Item(6, 19).Value2 = 1111
$ws.Cells.Item(6, 20).Value2 = 18

# Row 7
$ws.Cells.Item(7, 4).Value2 = 44328
$ws.Cells.Item(7, 12).Value2 = 'Especial'
$ws.Cells.Item(7, 13).Value2 = 250
$ws.Cells.Item(7, 14).Value2 = 20000
$ws.Cells.Item(7, 15).Value2 = 20000
$ws.Cells.Item(7, 16).Value2 = 20000
$ws.Cells.Item(7, 17).Value2 = '$/caja 18 kilos granel'
$ws.Cells.Item(7, 18).Value2 = 'Provincia de Limarí'
$ws.Cells.Item(7, 19).Value2 = 1111
$ws.Cells.Item(7, 20).Value2 = 18

# Row 8
$ws.Cells.Item(8, 4).Value2 = 44319
$ws.Cells.Item(8, 12).Value2 = 'Especial'
$ws.Cells.Item(8, 13).Value2 = 120
$ws.Cells.Item(8, 14).Value2 = 20000
$ws.Cells.Item(8, 15).Value2 = 20000
$ws.Cells.Item(8, 16).Value2 = 20000
$ws.Cells.Item(8, 17).Value2 = '$/caja 18 kilos granel'
$ws.Cells.Item(8, 18).Value2 = 'Provincia de Limarí'
$ws.Cells.Item(8, 19).Value2 = 1111
$ws.Cells.Item(8, 20).Value2 = 18

# Row 9
$ws.Cells.Item(9, 4).Value2 = 44348
$ws.Cells.Item(9, 12).Value2 = 'Especial'
$ws.Cells.Item(9, 13).Value2 = 200
$ws.Cells.Item(9, 14).Value2 = 20000
$ws.Cells.Item(9, 15).Value2 = 20000
$ws.Cells.Item(9, 16).Value2 = 20000
$ws.Cells.Item(9, 17).Value2 = '$/caja 18 kilos granel'
$ws.Cells.Item(9, 18).Value2 = 'Provincia de Limarí'
$ws.Cells.Item(9, 19).Value2 = 1111
$ws.Cells.Item(9, 20).Value2 = 18

# Row 11
$ws.Cells.Item(11, 4).Value2 = 44294
$ws.Cells.Item(11, 12).Value2 = 'Primera'
$ws.Cells.Item(11, 13).Value2 = 50
$ws.Cells.Item(11, 14).Value2 = 12000
$ws.Cells.Item(11, 15).Value2 = 12000
$ws.Cells.Item(11, 16).Value2 = 12000
$ws.Cells.Item(11, 17).Value2 = '$/caja 15 kilos granel'
$ws.Cells.Item(11, 18).Value2 = 'Región Metropolitana'
$ws.Cells.Item(11, 19).Value2 = 800
$ws.Cells.Item(11, 20).Value2 = 15

# Row 12
$ws.Cells.Item(12, 4).Value2 = 44692
$ws.Cells.Item(12, 12).Value2 = 'Especial'
$ws.Cells.Item(12, 13).Value2 = 150
$ws.Cells.Item(12, 14).Value2 = 17000
$ws.Cells.Item(12, 15).Value2 = 17000
$ws.Cells.Item(12, 16).Value2 = 17000
$ws.Cells.Item(12, 17).Value2 = '$/caja 18 kilos granel'
$ws.Cells.Item(12, 18).Value2 = 'Provincia de Limarí'
$ws.Cells.Item(12, 19).Value2 = 944
$ws.Cells.Item(12, 20).Value2 = 18

# Row 13
$ws.Cells.Item(13, 4).Value2 = 44354
$ws.Cells.Item(13, 12).Value2 = 'Primera'
$ws.Cells.Item(13, 13).Value2 = 100
$ws.Cells.Item(13, 14).Value2 = 18000
$ws.Cells.Item(13, 15).Value2 = 18000
$ws.Cells.Item(13, 16).Value2 = 18000
$ws.Cells.Item(13, 17).Value2 = '$/caja 18 kilos granel'
$ws.Cells.Item(13, 18).Value2 = 'Provincia de Limarí'
$ws.Cells.Item(13, 19).Value2 = 1000
$ws.Cells.Item(13, 20).Value2 = 18

# Row 14
$ws.Cells.Item(14, 4).Value2 = 44358
$ws.Cells.Item(14, 12).Value2 = 'Especial'
$ws.Cells.Item(14, 13).Value2 = 150
$ws.Cells.Item(14, 14).Value2 = 18000
$ws.Cells.Item(14, 15).Value2 = 18000
$ws.Cells.Item(14, 16).Value2 = 18000
$ws.Cells.Item(14, 17).Value2 = '$/caja 18 kilos granel'
$ws.Cells.Item(14, 18).Value2 = 'Provincia de Limarí'
$ws.Cells.Item(14, 19).Value2 = 1000
$ws.Cells.Item(14, 20).Value2 = 18

# Row 15
$ws.Cells.Item(15, 4).Value2 = 44358
$ws.Cells.Item(15, 12).Value2 = 'Primera'
$ws.Cells.Item(15, 13).Value2 = 100
$ws.Cells.Item(15, 14).Value2 = 17000
$ws.Cells.Item(15, 15).Value2 = 17000
$ws.Cells.Item(15, 16).Value2 = 17000
$ws.Cells.Item(15, 17).Value2 = '$/caja 18 kilos granel'
$ws.Cells.Item(15, 18).Value2 = 'Provincia de Limarí'
$ws.Cells.Item(15, 19).Value2 = 944
$ws.Cells.Item(15, 20).Value2 = 18

# Row 16
$ws.Cells.Item(16, 4).Value2 = 44355
$ws.Cells.Item(16, 12).Value2 = 'Especial'
$ws.Cells.Item(16, 13).Value2 = 50
$ws.Cells.Item(16, 14).Value2 = 18000
$ws.Cells.Item(16, 15).Value2 = 18000
$ws.Cells.Item(16, 16).Value2 = 18000
$ws.Cells.Item(16, 17).Value2 = '$/caja 18 kilos granel'
$ws.Cells.Item(16, 18).Value2 = 'Provincia de Limarí'
$ws.Cells.Item(16, 19).Value2 = 1000
$ws.Cells.Item(16, 20).Value2 = 18

# Row 17
$ws.Cells.Item(17, 4).Value2 = 44340
$ws.Cells.Item(17, 12).Value2 = 'Primera'
$ws.Cells.Item(17, 13).Value2 = 230
$ws.Cells.Item(17, 14).Value2 = 20000
$ws.Cells.Item(17, 15).Value2 = 20000
$ws.Cells.Item(17, 16).Value2 = 20000
$ws.Cells.Item(17, 17).Value2 = '$/caja 18 kilos granel'
$ws.Cells.Item(17, 18).Value2 = 'Provincia de Limarí'
$ws.Cells.Item(17, 19).Value2 = 1111
$ws.Cells.Item(17, 20).Value2 = 18

# Row 18
$ws.Cells.Item(18, 4).Value2 = 44326
$ws.Cells.Item(18, 12).Value2 = 'Especial'
$ws.Cells.Item(18, 13).Value2 = 300
$ws.Cells.Item(18, 14).Value2 = 20000
$ws.Cells.Item(18, 15).Value2 = 20000
$ws.Cells.Item(18, 16).Value2 = 20000
$ws.Cells.Item(18, 17).Value2 = '$/caja 18 kilos granel'
$ws.Cells.Item(18, 18).Value2 = 'Provincia de Limarí'
$ws.Cells.Item(18, 19).Value2 = 1111
$ws.Cells.Item(18, 20).Value2 = 18

# Row 19
$ws.Cells.Item(19, 4).Value2 = 44714
$ws.Cells.Item(19, 12).Value2 = 'Primera'
$ws.Cells.Item(19, 13).Value2 = 100
$ws.Cells.Item(19, 14).Value2 = 20000
$ws.Cells.Item(19, 15).Value2 = 20000
$ws.Cells.Item(19, 16).Value2 = 20000
$ws.Cells.Item(19, 17).Value2 = '$/caja 18 kilos granel'
$ws.Cells.Item(19, 18).Value2 = 'Provincia de Limarí'
$ws.Cells.Item(19, 19).Value2 = 1111
$ws.Cells.Item(19, 20).Value2 = 18
